$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 updates (columns K through S)
$ws.Range("K6").Value = 15.950976
$ws.Range("L6").Value = 17.992226
$ws.Range("M6").Value = 19.272346
$ws.Range("N6").Value = 20.648446
$ws.Range("O6").Value = 21.576324
$ws.Range("P6").Value = 22.964944
$ws.Range("Q6").Value = 23.928069
$ws.Range("R6").Value = 25.073389
$ws.Range("S6").Value = 27.099459

# Row 10 updates (columns K through S)
$ws.Range("K10").Value = 8.563734999999999
$ws.Range("L10").Value = 8.6057474
$ws.Range("M10").Value = 8.6057474
$ws.Range("N10").Value = 11.3516474
$ws.Range("O10").Value = 13.5204574
$ws.Range("P10").Value = 15.0479374
$ws.Range("Q10").Value = 16.2902874
$ws.Range("R10").Value = 17.7959774
$ws.Range("S10").Value = 19.2935374
